# MorganPatrone2006a__C_Stationarygenerator_alpha_non_zero.xlsx
# "expermits todos no convexos menos el 5to"
#
# Updates the numeric/expression inputs of the non-convex experiment
# generator workbook (leader/follower restrictions, the modified point,
# and the vec_bf / vec_BF / vec_alpha outputs) to the values of the new
# experiment run.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellRef, $text) {
    # Several of the target values are digit-strings (e.g. "-5.0", "0.34")
    # that Excel's normal Value assignment auto-converts to numbers. Forcing
    # the cell to Text format for the duration of the write keeps them as
    # literal text (matching the source workbook, where these are strings),
    # then restoring the "Normal" style leaves no stray formatting behind.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Restricciones_del_lider ---------------------------------------------
$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
$wsLider.Range("A2").Value = "4.5 - x"
Set-TextValue $wsLider "B2" "-5.0"
Set-TextValue $wsLider "D2" "0.34"
$wsLider.Range("A3").Value = "-4.5 + x"
Set-TextValue $wsLider "B3" "4.0"
Set-TextValue $wsLider "D3" "0.0"

# --- Restricciones_del_follower ------------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
$wsFollower.Range("A2").Value = "8.872427983539092 - 3.1687242798353905y"
Set-TextValue $wsFollower "B2" "-9.872427983539092"
Set-TextValue $wsFollower "D2" "0.09"
Set-TextValue $wsFollower "E2" "-5.6000000000000005"
Set-TextValue $wsFollower "F2" "-7.7"
$wsFollower.Range("A3").Value = "4.0040000000000004 - 1.4300000000000002y"
Set-TextValue $wsFollower "B3" "-5.0040000000000004"
Set-TextValue $wsFollower "D3" "0.82"
Set-TextValue $wsFollower "E3" "0.4"
Set-TextValue $wsFollower "F3" "0"

# --- Punto_modificado ------------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto "A2" "4.5"
Set-TextValue $wsPunto "B2" "2.8"

# --- Vector_bf / Vector_BF ------------------------------------------------
# NOTE: worksheet lookup by name is case-insensitive in this host, and the
# workbook has both a "Vector_bf" and a "Vector_BF" sheet, so those two are
# addressed by their (stable) tab position instead of by name.
$wsBf = $wb.Worksheets.Item(5)   # Vector_bf
Set-TextValue $wsBf "A2" "-3.0422148148148147"

$wsBF = $wb.Worksheets.Item(6)   # Vector_BF
Set-TextValue $wsBF "A2" "-1.09"
Set-TextValue $wsBF "A3" "-16.17285596707819"

# --- Vector_Alpha --------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 2.43
